$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37934
$ws1.Range("F4").Value = 642
$ws1.Range("F6").Value = 489
$ws1.Range("F7").Value = 376
$ws1.Range("F9").Value = 865
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 746
$ws1.Range("F12").Value = 581
$ws1.Range("F13").Value = 79
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 685
$ws1.Range("F18").Value = 486
$ws1.Range("F19").Value = 448
$ws1.Range("F21").Value = 98
$ws1.Range("F22").Value = 874
$ws1.Range("F23").Value = 2585
$ws1.Range("F24").Value = 1073
$ws1.Range("F26").Value = 116
$ws1.Range("F27").Value = 1176
$ws1.Range("F29").Value = 832
$ws1.Range("F30").Value = 75
$ws1.Range("F31").Value = 1179

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 441

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 666

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 666
$ws4.Range("F3").Value = 37934
$ws4.Range("F5").Value = 642
$ws4.Range("F7").Value = 489
$ws4.Range("F9").Value = 376
$ws4.Range("F11").Value = 441
$ws4.Range("F12").Value = 442
$ws4.Range("F13").Value = 338
$ws4.Range("F16").Value = 865
$ws4.Range("F17").Value = 106
$ws4.Range("F18").Value = 746
$ws4.Range("F19").Value = 581
$ws4.Range("F20").Value = 79
$ws4.Range("F26").Value = 35
$ws4.Range("F28").Value = 685
$ws4.Range("F30").Value = 486
$ws4.Range("F31").Value = 448
$ws4.Range("F33").Value = 98
$ws4.Range("F34").Value = 874
$ws4.Range("F35").Value = 2585
$ws4.Range("F36").Value = 1073
$ws4.Range("F38").Value = 116
$ws4.Range("F39").Value = 1176
$ws4.Range("F42").Value = 832
$ws4.Range("F43").Value = 75
$ws4.Range("F44").Value = 1179
